# Natmi following Dr Hou advice
# Rewrites the LR-pairs data rows (2-16) to include the new "FAPs" target-cluster
# rows and refreshed NATMI statistics for every Sending/Target cluster pairing.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> FAPs (App/Fpr2)
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "App"
$ws.Cells.Item(2, 3).Value = "Fpr2"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 126.9318136666667
$ws.Cells.Item(2, 8).Value = 380.795441
$ws.Cells.Item(2, 9).Value = 0.1973293860115714
$ws.Cells.Item(2, 10).Value = 0.1973293860115715
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 1.109174333333333
$ws.Cells.Item(2, 14).Value = 3.327523
$ws.Cells.Item(2, 15).Value = 0.06061832081580493
$ws.Cells.Item(2, 16).Value = 0.06061832081580493
$ws.Cells.Item(2, 17).Value = 140.7895098025159
$ws.Cells.Item(2, 18).Value = 1267.105588222643
$ws.Cells.Item(2, 19).Value = 0.01196177602763525
$ws.Cells.Item(2, 20).Value = 0.01196177602763525

# Row 3: ECs -> M1 (App/Fpr2)
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "App"
$ws.Cells.Item(3, 3).Value = "Fpr2"
$ws.Cells.Item(3, 4).Value = "M1"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 126.9318136666667
$ws.Cells.Item(3, 8).Value = 380.795441
$ws.Cells.Item(3, 9).Value = 0.1973293860115714
$ws.Cells.Item(3, 10).Value = 0.1973293860115715
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 4.402069333333333
$ws.Cells.Item(3, 14).Value = 13.206208
$ws.Cells.Item(3, 15).Value = 0.2405808023879173
$ws.Cells.Item(3, 16).Value = 0.2405808023879173
$ws.Cells.Item(3, 17).Value = 558.7626443664142
$ws.Cells.Item(3, 18).Value = 5028.863799297727
$ws.Cells.Item(3, 19).Value = 0.04747366202137891
$ws.Cells.Item(3, 20).Value = 0.04747366202137893

# Row 4: ECs -> M2 (App/Fpr2)
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "App"
$ws.Cells.Item(4, 3).Value = "Fpr2"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 126.9318136666667
$ws.Cells.Item(4, 8).Value = 380.795441
$ws.Cells.Item(4, 9).Value = 0.1973293860115714
$ws.Cells.Item(4, 10).Value = 0.1973293860115715
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 12.78643133333333
$ws.Cells.Item(4, 14).Value = 38.359294
$ws.Cells.Item(4, 15).Value = 0.6988008767962779
$ws.Cells.Item(4, 16).Value = 0.6988008767962779
$ws.Cells.Item(4, 17).Value = 1623.004919464295
$ws.Cells.Item(4, 18).Value = 14607.04427517865
$ws.Cells.Item(4, 19).Value = 0.1378939479625573
$ws.Cells.Item(4, 20).Value = 0.1378939479625573

# Row 5: FAPs -> FAPs (App/Fpr2)
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "App"
$ws.Cells.Item(5, 3).Value = "Fpr2"
$ws.Cells.Item(5, 4).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 152.3944216666667
$ws.Cells.Item(5, 8).Value = 457.183265
$ws.Cells.Item(5, 9).Value = 0.2369137948193439
$ws.Cells.Item(5, 10).Value = 0.2369137948193439
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 1.109174333333333
$ws.Cells.Item(5, 14).Value = 3.327523
$ws.Cells.Item(5, 15).Value = 0.06061832081580493
$ws.Cells.Item(5, 16).Value = 0.06061832081580493
$ws.Cells.Item(5, 17).Value = 169.0319810558439
$ws.Cells.Item(5, 18).Value = 1521.287829502595
$ws.Cells.Item(5, 19).Value = 0.01436131642004877
$ws.Cells.Item(5, 20).Value = 0.01436131642004877

# Row 6: FAPs -> M1 (App/Fpr2)
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "App"
$ws.Cells.Item(6, 3).Value = "Fpr2"
$ws.Cells.Item(6, 4).Value = "M1"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 152.3944216666667
$ws.Cells.Item(6, 8).Value = 457.183265
$ws.Cells.Item(6, 9).Value = 0.2369137948193439
$ws.Cells.Item(6, 10).Value = 0.2369137948193439
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 4.402069333333333
$ws.Cells.Item(6, 14).Value = 13.206208
$ws.Cells.Item(6, 15).Value = 0.2405808023879173
$ws.Cells.Item(6, 16).Value = 0.2405808023879173
$ws.Cells.Item(6, 17).Value = 670.8508101899022
$ws.Cells.Item(6, 18).Value = 6037.65729170912
$ws.Cells.Item(6, 19).Value = 0.05699691085440415
$ws.Cells.Item(6, 20).Value = 0.05699691085440416

# Row 7: FAPs -> M2 (App/Fpr2)
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "App"
$ws.Cells.Item(7, 3).Value = "Fpr2"
$ws.Cells.Item(7, 4).Value = "M2"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 152.3944216666667
$ws.Cells.Item(7, 8).Value = 457.183265
$ws.Cells.Item(7, 9).Value = 0.2369137948193439
$ws.Cells.Item(7, 10).Value = 0.2369137948193439
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 12.78643133333333
$ws.Cells.Item(7, 14).Value = 38.359294
$ws.Cells.Item(7, 15).Value = 0.6988008767962779
$ws.Cells.Item(7, 16).Value = 0.6988008767962779
$ws.Cells.Item(7, 17).Value = 1948.580808223879
$ws.Cells.Item(7, 18).Value = 17537.22727401491
$ws.Cells.Item(7, 19).Value = 0.165555567544891
$ws.Cells.Item(7, 20).Value = 0.165555567544891

# Row 8: M1 -> FAPs (App/Fpr2)
$ws.Cells.Item(8, 1).Value = "M1"
$ws.Cells.Item(8, 2).Value = "App"
$ws.Cells.Item(8, 3).Value = "Fpr2"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 206.573929
$ws.Cells.Item(8, 8).Value = 619.7217869999999
$ws.Cells.Item(8, 9).Value = 0.321141764212203
$ws.Cells.Item(8, 10).Value = 0.321141764212203
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 1.109174333333333
$ws.Cells.Item(8, 14).Value = 3.327523
$ws.Cells.Item(8, 15).Value = 0.06061832081580493
$ws.Cells.Item(8, 16).Value = 0.06061832081580493
$ws.Cells.Item(8, 17).Value = 229.1264999826223
$ws.Cells.Item(8, 18).Value = 2062.138499843601
$ws.Cells.Item(8, 19).Value = 0.01946707449036891
$ws.Cells.Item(8, 20).Value = 0.01946707449036891

# Row 9: M1 -> M1 (App/Fpr2)
$ws.Cells.Item(9, 1).Value = "M1"
$ws.Cells.Item(9, 2).Value = "App"
$ws.Cells.Item(9, 3).Value = "Fpr2"
$ws.Cells.Item(9, 4).Value = "M1"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 206.573929
$ws.Cells.Item(9, 8).Value = 619.7217869999999
$ws.Cells.Item(9, 9).Value = 0.321141764212203
$ws.Cells.Item(9, 10).Value = 0.321141764212203
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 4.402069333333333
$ws.Cells.Item(9, 14).Value = 13.206208
$ws.Cells.Item(9, 15).Value = 0.2405808023879173
$ws.Cells.Item(9, 16).Value = 0.2405808023879173
$ws.Cells.Item(9, 17).Value = 909.3527579170773
$ws.Cells.Item(9, 18).Value = 8184.174821253696
$ws.Cells.Item(9, 19).Value = 0.07726054331444314
$ws.Cells.Item(9, 20).Value = 0.07726054331444314

# Row 10: M1 -> M2 (App/Fpr2)
$ws.Cells.Item(10, 1).Value = "M1"
$ws.Cells.Item(10, 2).Value = "App"
$ws.Cells.Item(10, 3).Value = "Fpr2"
$ws.Cells.Item(10, 4).Value = "M2"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 206.573929
$ws.Cells.Item(10, 8).Value = 619.7217869999999
$ws.Cells.Item(10, 9).Value = 0.321141764212203
$ws.Cells.Item(10, 10).Value = 0.321141764212203
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 12.78643133333333
$ws.Cells.Item(10, 14).Value = 38.359294
$ws.Cells.Item(10, 15).Value = 0.6988008767962779
$ws.Cells.Item(10, 16).Value = 0.6988008767962779
$ws.Cells.Item(10, 17).Value = 2641.343358415375
$ws.Cells.Item(10, 18).Value = 23772.09022573838
$ws.Cells.Item(10, 19).Value = 0.224414146407391
$ws.Cells.Item(10, 20).Value = 0.224414146407391

# Row 11: M2 -> FAPs (App/Fpr2)
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "App"
$ws.Cells.Item(11, 3).Value = "Fpr2"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 141.7744496666667
$ws.Cells.Item(11, 8).Value = 425.323349
$ws.Cells.Item(11, 9).Value = 0.2204038869114384
$ws.Cells.Item(11, 10).Value = 0.2204038869114385
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 1.109174333333333
$ws.Cells.Item(11, 14).Value = 3.327523
$ws.Cells.Item(11, 15).Value = 0.06061832081580493
$ws.Cells.Item(11, 16).Value = 0.06061832081580493
$ws.Cells.Item(11, 17).Value = 157.2525806927252
$ws.Cells.Item(11, 18).Value = 1415.273226234527
$ws.Cells.Item(11, 19).Value = 0.01336051352584797
$ws.Cells.Item(11, 20).Value = 0.01336051352584797

# Row 12: M2 -> M1 (App/Fpr2)
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "App"
$ws.Cells.Item(12, 3).Value = "Fpr2"
$ws.Cells.Item(12, 4).Value = "M1"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 141.7744496666667
$ws.Cells.Item(12, 8).Value = 425.323349
$ws.Cells.Item(12, 9).Value = 0.2204038869114384
$ws.Cells.Item(12, 10).Value = 0.2204038869114385
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 4.402069333333333
$ws.Cells.Item(12, 14).Value = 13.206208
$ws.Cells.Item(12, 15).Value = 0.2405808023879173
$ws.Cells.Item(12, 16).Value = 0.2405808023879173
$ws.Cells.Item(12, 17).Value = 624.1009571278436
$ws.Cells.Item(12, 18).Value = 5616.908614150592
$ws.Cells.Item(12, 19).Value = 0.05302494396256964
$ws.Cells.Item(12, 20).Value = 0.05302494396256965

# Row 13: M2 -> M2 (App/Fpr2)
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "App"
$ws.Cells.Item(13, 3).Value = "Fpr2"
$ws.Cells.Item(13, 4).Value = "M2"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 141.7744496666667
$ws.Cells.Item(13, 8).Value = 425.323349
$ws.Cells.Item(13, 9).Value = 0.2204038869114384
$ws.Cells.Item(13, 10).Value = 0.2204038869114385
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 12.78643133333333
$ws.Cells.Item(13, 14).Value = 38.359294
$ws.Cells.Item(13, 15).Value = 0.6988008767962779
$ws.Cells.Item(13, 16).Value = 0.6988008767962779
$ws.Cells.Item(13, 17).Value = 1812.789265483956
$ws.Cells.Item(13, 18).Value = 16315.10338935561
$ws.Cells.Item(13, 19).Value = 0.1540184294230209
$ws.Cells.Item(13, 20).Value = 0.1540184294230209

# Row 14: sCs -> FAPs (App/Fpr2)
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "App"
$ws.Cells.Item(14, 3).Value = "Fpr2"
$ws.Cells.Item(14, 4).Value = "FAPs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 15.57379533333333
$ws.Cells.Item(14, 8).Value = 46.721386
$ws.Cells.Item(14, 9).Value = 0.02421116804544314
$ws.Cells.Item(14, 10).Value = 0.02421116804544315
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 1.109174333333333
$ws.Cells.Item(14, 14).Value = 3.327523
$ws.Cells.Item(14, 15).Value = 0.06061832081580493
$ws.Cells.Item(14, 16).Value = 0.06061832081580493
$ws.Cells.Item(14, 17).Value = 17.27405405631977
$ws.Cells.Item(14, 18).Value = 155.466486506878
$ws.Cells.Item(14, 19).Value = 0.001467640351904037
$ws.Cells.Item(14, 20).Value = 0.001467640351904037

# Row 15: sCs -> M1 (App/Fpr2)
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "App"
$ws.Cells.Item(15, 3).Value = "Fpr2"
$ws.Cells.Item(15, 4).Value = "M1"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 15.57379533333333
$ws.Cells.Item(15, 8).Value = 46.721386
$ws.Cells.Item(15, 9).Value = 0.02421116804544314
$ws.Cells.Item(15, 10).Value = 0.02421116804544315
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 4.402069333333333
$ws.Cells.Item(15, 14).Value = 13.206208
$ws.Cells.Item(15, 15).Value = 0.2405808023879173
$ws.Cells.Item(15, 16).Value = 0.2405808023879173
$ws.Cells.Item(15, 17).Value = 68.55692684047644
$ws.Cells.Item(15, 18).Value = 617.0123415642879
$ws.Cells.Item(15, 19).Value = 0.005824742235121414
$ws.Cells.Item(15, 20).Value = 0.005824742235121415

# Row 16: sCs -> M2 (App/Fpr2)
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "App"
$ws.Cells.Item(16, 3).Value = "Fpr2"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 15.57379533333333
$ws.Cells.Item(16, 8).Value = 46.721386
$ws.Cells.Item(16, 9).Value = 0.02421116804544314
$ws.Cells.Item(16, 10).Value = 0.02421116804544315
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 12.78643133333333
$ws.Cells.Item(16, 14).Value = 38.359294
$ws.Cells.Item(16, 15).Value = 0.6988008767962779
$ws.Cells.Item(16, 16).Value = 0.6988008767962779
$ws.Cells.Item(16, 17).Value = 199.1332646290537
$ws.Cells.Item(16, 18).Value = 1792.199381661484
$ws.Cells.Item(16, 19).Value = 0.01691878545841769
$ws.Cells.Item(16, 20).Value = 0.0169187854584177
